$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename label strings in column A (these are shared-string backed labels).
# skala1item1 -> skala1_item1
$ws.Range("A6:A9").Value = "skala1_item1"

# skala1item2 -> skala1_item2
$ws.Range("A10:A13").Value = "skala1_item2"

# skala1item3 -> skala1_item3
$ws.Range("A14:A17").Value = "skala1_item3"

# pv_kat1 -> pvkat_1
$ws.Range("A18:A22").Value = "pvkat_1"

# pv_kat2 -> pvkat_2
$ws.Range("A23:A27").Value = "pvkat_2"

# pv_kat3 -> pvkat_3
$ws.Range("A28:A32").Value = "pvkat_3"

# pv_kat4 -> pvkat_4
$ws.Range("A33:A37").Value = "pvkat_4"

# pv_kat5 -> pvkat_5
$ws.Range("A38:A42").Value = "pvkat_5"
